# Swap the data of row 3 and row 4 (the two species records got
# reordered/re-matched), refreshing the Ost/Nord (Q/R) coordinates to the
# new rounded values that came with the swap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 gets what used to be row 4's data ------------------------------
$ws.Range("A3").Value  = 112144586
$ws.Range("B3").Value  = 44322
$ws.Range("E3").Value  = 102366
$ws.Range("F3").Value  = "Ängsmetallvinge"
$ws.Range("G3").Value  = "Adscita statices"
$ws.Range("I3").Value  = ""
$ws.Range("J3").Value  = ""
$ws.Range("K3").Value  = ""
$ws.Range("L3").Value  = ""
$ws.Range("M3").Value  = ""
$ws.Range("Q3").Value  = 442995
$ws.Range("R3").Value  = 6204827
$ws.Range("Z3").Value  = ""
$ws.Range("AB3").Value = ""
$ws.Range("AC3").Value = ""
$ws.Range("AO3").Value = ""
$ws.Range("AQ3").Value = ""
$ws.Range("AR3").Value = ""

# --- Row 4 gets what used to be row 3's data ------------------------------
$ws.Range("A4").Value  = 112144581
$ws.Range("B4").Value  = 42594
$ws.Range("E4").Value  = 101260
$ws.Range("F4").Value  = "Svartfläckig blåvinge"
$ws.Range("G4").Value  = "Phengaris arion"
$ws.Range("I4").Value  = "1"
$ws.Range("J4").Value  = "ex."
$ws.Range("K4").Value  = "imago/adult"
$ws.Range("L4").Value  = "hona"
$ws.Range("M4").Value  = "vilande"
$ws.Range("Q4").Value  = 442972
$ws.Range("R4").Value  = 6204767
$ws.Range("Z4").Value  = ""
$ws.Range("AB4").Value = ""
$ws.Range("AC4").Value = "lufthåvning"
$ws.Range("AO4").Value = "på grässtrå"
$ws.Range("AQ4").Value = "Nils Otto Nilsson"
$ws.Range("AR4").Value = "NON 04616"
